# mTBIconditions4.xlsx — fix: condition files, psychopy, xmls and fake data
#
# The original "D" column (opacity, always 1) is removed entirely and the
# "image"/"q_or_r" columns get their values normalised:
#   up.jpg / down.jpg          -> images/up.jpg / images/down.jpg
#   "Question?" and "Rest"     -> "Rate your ability to control your brain"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- drop the now-unused opacity column (D) -------------------------------
$ws.Columns.Item(4).Delete()

# --- point the image file names at the images/ folder ---------------------
$ws.Range("B2").Value = "images/up.jpg"
$ws.Range("B5").Value = "images/up.jpg"
$ws.Range("B6").Value = "images/up.jpg"

# --- collapse the old "Question?"/"Rest" prompts into one ----------------
$ws.Range("C2").Value = "Rate your ability to control your brain"
$ws.Range("C3").Value = "Rate your ability to control your brain"
$ws.Range("C4").Value = "Rate your ability to control your brain"
$ws.Range("C5").Value = "Rate your ability to control your brain"
$ws.Range("C6").Value = "Rate your ability to control your brain"
$ws.Range("C7").Value = "Rate your ability to control your brain"

$ws.Range("B3").Value = "images/down.jpg"
$ws.Range("B4").Value = "images/down.jpg"
$ws.Range("B7").Value = "images/down.jpg"

# --- widen the three remaining columns to fit the new content -------------
$ws.Columns.Item(1).ColumnWidth = 11.6313725490196
$ws.Columns.Item(2).ColumnWidth = 15.121568627451
$ws.Columns.Item(3).ColumnWidth = 31.321568627451

# --- move the active selection from the old E16 to D1 ---------------------
$null = $ws.Range("D1").Select()
